$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.608.55"
$ws.Range("E2").Value = "  +3.95%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.417.07"
$ws.Range("E3").Value = "  +3.28%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.08"
$ws.Range("E5").Value = "  +2.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.58"
$ws.Range("E6").Value = "  +3.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  +3.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.411.34"
$ws.Range("E8").Value = "  +3.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +14.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.632"
$ws.Range("E11").Value = "  +3.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.08"
$ws.Range("E12").Value = "  +4.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000282"
$ws.Range("E13").Value = "  +7.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.15"
$ws.Range("E14").Value = "  +3.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.966.73"
$ws.Range("E15").Value = "  +3.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.38"
$ws.Range("E16").Value = "  +4.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.420.79"
$ws.Range("E17").Value = "  +3.20%  "

$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "65.597.49"
$ws.Range("E19").Value = "  +4.05%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.96"
$ws.Range("E20").Value = "  +3.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").Value = "  +3.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.71"
$ws.Range("E22").Value = "  +16.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.02"
$ws.Range("E23").Value = "  +18.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.15"
$ws.Range("E24").Value = "  +3.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.57"
$ws.Range("E25").Value = "  +5.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.65"
$ws.Range("E26").Value = "  +4.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.90"
$ws.Range("E27").Value = "  +3.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.90"
$ws.Range("E28").Value = "  +7.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.89"
$ws.Range("E29").Value = "  +4.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.82"
$ws.Range("E30").Value = "  +6.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.70"
$ws.Range("E31").Value = "  +5.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.57"
$ws.Range("E32").Value = "  +2.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "586.49"
$ws.Range("E33").Value = "  +3.26%  "

$ws.Range("E34").Value = "  +3.71%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "61.12"
$ws.Range("E35").Value = "  +6.73%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.141"
$ws.Range("E37").Value = "  -3.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.04"
$ws.Range("E38").Value = "  +2.96%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.50"
$ws.Range("E39").Value = "  +2.38%  "

$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0760"
$ws.Range("E40").Value = "  +3.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.375"
$ws.Range("E41").Value = "  +3.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.108.36"
$ws.Range("E42").Value = "  -1.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.89"
$ws.Range("E44").Value = "  +3.69%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0417"
$ws.Range("E45").Value = "  +4.44%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.51"
$ws.Range("E46").Value = "  +4.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.23"
$ws.Range("E47").Value = "  +1.87%  "

$ws.Range("E48").Value = "  +5.84%  "

$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.39"
$ws.Range("E50").Value = "  +5.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.42"
$ws.Range("E51").Value = "  +2.88%  "
